# Add a missing product row ("SOSRO TEH KOTAK 4+2S") into the BEV03S
# listing on Sheet1. The row belongs right after "SOSRO TEH KOTAK B250"
# (row 109) and before "NESCAFE CPUCCINO 220" (old row 110), continuing
# the D/E numbering sequence (category 11, sequence 6) that was
# previously skipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 110, shifting everything below down by one.
$ws.Rows(110).Insert()

# Copy the formatting (style/border) of the row above so the new row
# matches the rest of the table (thin-bordered "Normal" style) instead
# of picking up the blank default style.
$ws.Range("A109:F109").Copy()
$ws.Range("A110:F110").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row. A/D/E look numeric ("20001119", "11", "6") but
# the whole sheet stores them as text, so route them through a
# temporary text formula and then flatten back to a static value -
# this avoids Excel's automatic number coercion while keeping the
# cell's number format/style untouched.
$ws.Range("A110").Formula = '="20001119"'
$ws.Range("B110").Value = "SOSRO TEH KOTAK 4+2S"
$ws.Range("C110").Value = "BEV03S"
$ws.Range("D110").Formula = '="11"'
$ws.Range("E110").Formula = '="6"'
$ws.Range("F110").Value = "RT,(E-1B)"

# Convert the formula cells to plain static text values.
$ws.Range("A110:F110").Copy()
$ws.Range("A110:F110").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
